# Add the "Iphone 15" row to the iPhone release-dates table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fonts/fills/borders/number formats) of the last
# data row (row 20) down onto the new row 21, then overwrite values.
$ws.Range("A20:D20").Copy()
$ws.Range("A21:D21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(21, 1).Value = "Iphone 15"
$ws.Cells.Item(21, 2).Value = 2023
$ws.Cells.Item(21, 3).Value = 45557
$ws.Cells.Item(21, 4).Value = 45191

$ws.Rows(21).RowHeight = 27

# Row 22 mirrors the thin trailing formatting row under the table, only
# column C keeps the same formatting as the rest of the date column.
$ws.Range("C20").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows(22).RowHeight = 26

$ws.Range("D24").Select()
